# Update "想去人数" (want-to-go count) values in column F for rows 2,3,5,6,8
# on both the "展览" and "全部类型" worksheets, matching the data refresh
# captured by the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    2 = 2281
    3 = 1737
    5 = 1094
    6 = 860
    8 = 5850
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
